$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    3  = 5.951300000000004
    14 = 5.691199999999998
    21 = 9.551300000000007
    23 = 9.008900000000002
    25 = 6.061400000000001
    26 = 4.973300000000006
    29 = 5.106900000000005
    53 = 5.422799999999999
    57 = 4.881299999999997
    59 = 4.840299999999996
    69 = 5.339599999999995
    79 = 9.051700000000004
    83 = 5.357699999999997
    91 = 5.233200000000002
    93 = 5.828399999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}

$wb.Save()
